# Append the 2026-02-11 (Excel serial 46064) daily block to Daily_Data.
# Mirrors the 22-row-per-day pattern already present in the sheet:
# one "Registered" + one "Eligible" row per depository, in the same
# depository order used for every earlier day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Data")

$dateSerial = 46064

$newRows = @(
    @("ASAHI DEPOSITORY LLC Registered", 0, 0),
    @("ASAHI DEPOSITORY LLC Eligible", 0, 0),
    @("BRINK'S, INC. Registered", 73354.783, 73354.783),
    @("BRINK'S, INC. Eligible", 85821.84699999999, 85821.84699999999),
    @("CNT DEPOSITORY, INC. Registered", 1246.06, 1246.06),
    @("CNT DEPOSITORY, INC. Eligible", 0, 0),
    @("DELAWARE DEPOSITORY Registered", 1633.941, 1633.941),
    @("DELAWARE DEPOSITORY Eligible", 18459.584, 18459.584),
    @("HSBC BANK, USA Registered", 1394.758, 1394.758),
    @("HSBC BANK, USA Eligible", 9281.978999999999, 9281.978999999999),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 2395.448, 2395.448),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 0, 0),
    @("JP MORGAN CHASE BANK NA Registered", 114061.421, 114061.421),
    @("JP MORGAN CHASE BANK NA Eligible", 76408.66899999999, 76408.66899999999),
    @("LOOMIS INTERNATIONAL (US) LLC Registered", 61157.444, 61157.444),
    @("LOOMIS INTERNATIONAL (US) LLC Eligible", 71594.18700000001, 71594.18700000001),
    @("MALCA-AMIT USA, LLC Registered", 395.145, 395.145),
    @("MALCA-AMIT USA, LLC Eligible", 0, 0),
    @("MANFRA, TORDELLA & BROOKES, LLC Registered", 49920.248, 49920.248),
    @("MANFRA, TORDELLA & BROOKES, LLC Eligible", 2104.855, 2104.855),
    @("STONEX PRECIOUS METALS LLC Registered", 14122.765, 14122.765),
    @("STONEX PRECIOUS METALS LLC Eligible", 16.075, 16.075)
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$dateNumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $lastRow + 1 + $i
    $rowData = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = $dateNumberFormat
    $dateCell.Value = $dateSerial

    $ws.Cells.Item($r, 2).Value = $rowData[0]
    $ws.Cells.Item($r, 3).Value = $rowData[1]
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 7).Value = 0
    $ws.Cells.Item($r, 8).Value = $rowData[2]
}

Write-Output "Appended $($newRows.Count) rows starting at row $($lastRow + 1)"
